# The data table (Sheet1, A1:I49) lists monthly index values, one row per
# month, grouped in 12-row blocks per year (2014, 2015, 2016, 2017) starting
# at row 2. Within each year's block the rows were re-ordered so the last
# three months (Oct, Nov, Dec) come first, followed by Jan..Sep, i.e. each
# 12-row block is rotated so rows 10-12 move to the top (rows 1-9 follow).
#
# This reads each year block as a 2-D array via Range.Value2, builds the
# rotated block in memory, and writes it back with Range.Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$rowsPerYear = 12
$years = 4
$lastCol = 9

for ($y = 0; $y -lt $years; $y++) {
    $startRow = $firstDataRow + ($y * $rowsPerYear)
    $endRow = $startRow + $rowsPerYear - 1

    $rangeAddr = "A" + $startRow + ":I" + $endRow
    $block = $ws.Range($rangeAddr).Value2

    # Build the rotated block: original rows 10,11,12 (Oct,Nov,Dec) first,
    # then original rows 1..9 (Jan..Sep).
    # Note: $block (from Range.Value2) is 1-based, but a freshly
    # New-Object'd managed array is 0-based.
    $rotated = New-Object 'object[,]' $rowsPerYear, $lastCol

    $srcOrder = @(10, 11, 12, 1, 2, 3, 4, 5, 6, 7, 8, 9)
    for ($destRowIdx = 0; $destRowIdx -lt $rowsPerYear; $destRowIdx++) {
        $srcRow = $srcOrder[$destRowIdx]
        for ($col = 1; $col -le $lastCol; $col++) {
            $rotated[$destRowIdx, ($col - 1)] = $block[$srcRow, $col]
        }
    }

    $ws.Range($rangeAddr).Value = $rotated
}
